# REVER_DailyTracker_BALAJI.xlsx - "Add files via upload" update
#
# This adds/updates Daily-Tracker rows for the APR-2021 sheet:
#  - Row 7 (No=6): Application/Task switched to a Holiday / Election entry
#  - Rows 20-24 (No=19..23): Application + Task comments filled in for the
#    Mujistore validation-message work
#  - Row 27 (No=26): Application + Task comment filled in for the Suma-san
#    resource-file work
#  - The existing "sumasen" typo in row 17's comment is corrected to "sumasan"
#  - The last worksheet selection moves from D17 to D27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing typo: "sumasen" -> "sumasan" (edits shared string in place) ---
$ws.Range("D17").Value = "Validation messages given by sumasan task going on"

# --- Rows 20-24: Mujistore application + task comments -------------------
# (Order matters: new shared strings are appended in the order they are
# first written, and must land at indices 52-58 to match the target file.)
$ws.Range("C20").Value = "Mujistore "
$ws.Range("D20").Value = "Validation messages and success messages given by sumasan task going on"

$ws.Range("C21").Value = "Mujistore "
$ws.Range("D21").Value = "Nearly 185 validation messages changes done in English language and Japanese languages going on"

$ws.Range("C22").Value = "Mujistore "
$ws.Range("D22").Value = "Nirmal san issues going on and testing going on and until completed task files sent to deployment…"

$ws.Range("C23").Value = "Mujistore "
$ws.Range("D23").Value = "Nirmal san issues done and testing going on"

$ws.Range("C24").Value = "Mujistore "
$ws.Range("D24").Value = "Testing done and sent to deployment"

# --- Row 7: Holiday / Election --------------------------------------------
$ws.Range("C7").Value = "Holiday"
$ws.Range("D7").Value = "Election"

# --- Row 27: Mujistore application + task comment -------------------------
$ws.Range("C27").Value = "Mujistore "
$ws.Range("D27").Value = "Suma san given Resource file changes done and sent to deployment and file upload development gather req and analyse the scenario."

# --- Update the active selection to D27 (matches the saved view state) ----
$ws.Range("D27").Select()
